$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 0.0008333333333333334
$ws.Range("K2").Value = 5835
$ws.Range("L2").Value = 0.01167
